# The deck's slide-master theme ("Integral" / "Red Violet" colour scheme,
# ppt/theme/theme1.xml) is re-coloured to the stock "Office Theme" colour
# scheme (ppt/theme/theme2.xml, used today only by the notes master).
#
# Only the 12 theme colour slots (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink)
# actually differ between the two themes in this deck - fonts and format
# scheme are identical - so re-painting the slide master's ThemeColorScheme
# reproduces the intended visual result.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Index -> (slot name, new RGB as 0xBBGGRR long, matching VBA's RGB() encoding)
# 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
# 9=accent5 10=accent6 11=hlink 12=folHlink
$colors.Item(1).RGB = 0          # dk1      000000
$colors.Item(2).RGB = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB = 6968388    # dk2      44546A
$colors.Item(4).RGB = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB = 3243501    # accent2  ED7D31
$colors.Item(7).RGB = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB = 49407      # accent4  FFC000
$colors.Item(9).RGB = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
